# Apply Price (D) and Volume(1h) (E) updates for existing rows, plus the
# newly-inserted "WrappedliquidstakedEther2.0" row that shifts rows 24-51
# down by one (dropping the prior last row, WOONetwork, off the 51-row table).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = "29.271.78"
$ws.Cells.Item(2, 5).Value = "  +0.12%  "

$ws.Cells.Item(3, 4).Value = "1.900.05"
$ws.Cells.Item(3, 5).Value = "  -0.34%  "

$ws.Cells.Item(4, 4).NumberFormat = "@"
$ws.Cells.Item(4, 4).Value = "1.001"
$ws.Cells.Item(4, 5).Value = "  -0.23%  "

$ws.Cells.Item(5, 4).NumberFormat = "@"
$ws.Cells.Item(5, 4).Value = "326.18"
$ws.Cells.Item(5, 5).Value = "  -0.54%  "

$ws.Cells.Item(6, 5).Value = "  -0.20%  "

$ws.Cells.Item(7, 4).NumberFormat = "@"
$ws.Cells.Item(7, 4).Value = "0.4634"
$ws.Cells.Item(7, 5).Value = "  -0.56%  "

$ws.Cells.Item(8, 5).Value = "  -0.98%  "

$ws.Cells.Item(9, 4).NumberFormat = "@"
$ws.Cells.Item(9, 4).Value = "0.07885"
$ws.Cells.Item(9, 5).Value = "  -1.06%  "

$ws.Cells.Item(10, 4).NumberFormat = "@"
$ws.Cells.Item(10, 4).Value = "0.9894"
$ws.Cells.Item(10, 5).Value = "  -1.37%  "

$ws.Cells.Item(11, 4).NumberFormat = "@"
$ws.Cells.Item(11, 4).Value = "21.86"

$ws.Cells.Item(12, 4).Value = "1.915.81"
$ws.Cells.Item(12, 5).Value = "  +2.29%  "

$ws.Cells.Item(13, 4).NumberFormat = "@"
$ws.Cells.Item(13, 4).Value = "7.068"
$ws.Cells.Item(13, 5).Value = "  -0.97%  "

$ws.Cells.Item(14, 4).NumberFormat = "@"
$ws.Cells.Item(14, 4).Value = "5.738"
$ws.Cells.Item(14, 5).Value = "  -0.90%  "

$ws.Cells.Item(15, 4).NumberFormat = "@"
$ws.Cells.Item(15, 4).Value = "0.06980"
$ws.Cells.Item(15, 5).Value = "  +0.32%  "

$ws.Cells.Item(16, 4).NumberFormat = "@"
$ws.Cells.Item(16, 4).Value = "88.33"
$ws.Cells.Item(16, 5).Value = "  -0.69%  "

$ws.Cells.Item(17, 4).NumberFormat = "@"
$ws.Cells.Item(17, 4).Value = "1.002"
$ws.Cells.Item(17, 5).Value = "  -0.15%  "

$ws.Cells.Item(18, 4).NumberFormat = "@"
$ws.Cells.Item(18, 4).Value = "0.000009975"
$ws.Cells.Item(18, 5).Value = "  -1.54%  "

$ws.Cells.Item(19, 4).NumberFormat = "@"
$ws.Cells.Item(19, 4).Value = "17.07"
$ws.Cells.Item(19, 5).Value = "  -0.80%  "

$ws.Cells.Item(20, 5).Value = "  -0.10%  "

$ws.Cells.Item(21, 4).Value = "29.273.86"
$ws.Cells.Item(21, 5).Value = "  +0.10%  "

$ws.Cells.Item(22, 4).NumberFormat = "@"
$ws.Cells.Item(22, 4).Value = "5.299"
$ws.Cells.Item(22, 5).Value = "  -1.22%  "

$ws.Cells.Item(23, 5).Value = "  -0.19%  "

$ws.Cells.Item(24, 2).Value = "WrappedliquidstakedEther2.0"
$ws.Cells.Item(24, 3).Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Cells.Item(24, 4).Value = "2.140.37"
$ws.Cells.Item(24, 5).Value = "  +1.48%  "

$ws.Cells.Item(25, 2).Value = "Toncoin"
$ws.Cells.Item(25, 3).Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Cells.Item(25, 4).NumberFormat = "@"
$ws.Cells.Item(25, 4).Value = "2.100"
$ws.Cells.Item(25, 5).Value = "  +2.02%  "

$ws.Cells.Item(26, 2).Value = "Monero"
$ws.Cells.Item(26, 3).Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Cells.Item(26, 4).NumberFormat = "@"
$ws.Cells.Item(26, 4).Value = "156.07"
$ws.Cells.Item(26, 5).Value = "  -0.47%  "

$ws.Cells.Item(27, 2).Value = "EthereumClassic"
$ws.Cells.Item(27, 3).Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Cells.Item(27, 4).NumberFormat = "@"
$ws.Cells.Item(27, 4).Value = "19.39"
$ws.Cells.Item(27, 5).Value = "  -0.87%  "

$ws.Cells.Item(28, 2).Value = "InternetComputer(DFINITY)"
$ws.Cells.Item(28, 3).Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Cells.Item(28, 4).NumberFormat = "@"
$ws.Cells.Item(28, 4).Value = "6.030"
$ws.Cells.Item(28, 5).Value = "  +2.11%  "

$ws.Cells.Item(29, 2).Value = "BitcoinCash"
$ws.Cells.Item(29, 3).Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Cells.Item(29, 4).NumberFormat = "@"
$ws.Cells.Item(29, 4).Value = "118.33"
$ws.Cells.Item(29, 5).Value = "  -0.98%  "

$ws.Cells.Item(30, 2).Value = "LidoDAOToken"
$ws.Cells.Item(30, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(30, 4).NumberFormat = "@"
$ws.Cells.Item(30, 4).Value = "1.892"
$ws.Cells.Item(30, 5).Value = "  -5.72%  "

$ws.Cells.Item(31, 2).Value = "Stellar"
$ws.Cells.Item(31, 3).Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Cells.Item(31, 4).NumberFormat = "@"
$ws.Cells.Item(31, 4).Value = "0.09364"
$ws.Cells.Item(31, 5).Value = "  -0.83%  "

$ws.Cells.Item(32, 2).Value = "ImmutableX"
$ws.Cells.Item(32, 3).Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Cells.Item(32, 4).NumberFormat = "@"
$ws.Cells.Item(32, 4).Value = "0.9042"
$ws.Cells.Item(32, 5).Value = "  -2.35%  "

$ws.Cells.Item(33, 2).Value = "Filecoin"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Cells.Item(33, 4).NumberFormat = "@"
$ws.Cells.Item(33, 4).Value = "5.258"
$ws.Cells.Item(33, 5).Value = "  -1.83%  "

$ws.Cells.Item(34, 2).Value = "ARBITRUM"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Cells.Item(34, 4).NumberFormat = "@"
$ws.Cells.Item(34, 4).Value = "1.324"
$ws.Cells.Item(34, 5).Value = "  -1.90%  "

$ws.Cells.Item(35, 2).Value = "HuobiToken"
$ws.Cells.Item(35, 3).Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Cells.Item(35, 4).NumberFormat = "@"
$ws.Cells.Item(35, 4).Value = "3.209"
$ws.Cells.Item(35, 5).Value = "  -1.53%  "

$ws.Cells.Item(36, 2).Value = "TrustWalletToken"
$ws.Cells.Item(36, 3).Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Cells.Item(36, 4).NumberFormat = "@"
$ws.Cells.Item(36, 4).Value = "1.191"
$ws.Cells.Item(36, 5).Value = "  +1.79%  "

$ws.Cells.Item(37, 2).Value = "Hedera"
$ws.Cells.Item(37, 3).Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Cells.Item(37, 4).NumberFormat = "@"
$ws.Cells.Item(37, 4).Value = "0.05782"
$ws.Cells.Item(37, 5).Value = "  -1.37%  "

$ws.Cells.Item(38, 2).Value = "VeChain"
$ws.Cells.Item(38, 3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(38, 4).NumberFormat = "@"
$ws.Cells.Item(38, 4).Value = "0.02087"
$ws.Cells.Item(38, 5).Value = "  -1.11%  "

$ws.Cells.Item(39, 2).Value = "Frax"
$ws.Cells.Item(39, 3).Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Cells.Item(39, 4).NumberFormat = "@"
$ws.Cells.Item(39, 4).Value = "1.000"
$ws.Cells.Item(39, 5).Value = "  -0.25%  "

$ws.Cells.Item(40, 2).Value = "FraxShare"
$ws.Cells.Item(40, 3).Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Cells.Item(40, 4).NumberFormat = "@"
$ws.Cells.Item(40, 4).Value = "7.721"
$ws.Cells.Item(40, 5).Value = "  -3.39%  "

$ws.Cells.Item(41, 2).Value = "TheSandbox"
$ws.Cells.Item(41, 3).Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Cells.Item(41, 4).NumberFormat = "@"
$ws.Cells.Item(41, 4).Value = "0.5711"
$ws.Cells.Item(41, 5).Value = "  -0.90%  "

$ws.Cells.Item(42, 2).Value = "Algorand"
$ws.Cells.Item(42, 3).Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Cells.Item(42, 4).NumberFormat = "@"
$ws.Cells.Item(42, 4).Value = "0.1788"
$ws.Cells.Item(42, 5).Value = "  -1.50%  "

$ws.Cells.Item(43, 2).Value = "Aptos"
$ws.Cells.Item(43, 3).Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Cells.Item(43, 4).NumberFormat = "@"
$ws.Cells.Item(43, 4).Value = "9.720"
$ws.Cells.Item(43, 5).Value = "  -3.24%  "

$ws.Cells.Item(44, 2).Value = "EnergySwap"
$ws.Cells.Item(44, 3).Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Cells.Item(44, 4).NumberFormat = "@"
$ws.Cells.Item(44, 4).Value = "11.95"
$ws.Cells.Item(44, 5).Value = "  -0.54%  "

$ws.Cells.Item(45, 2).Value = "Decentraland"
$ws.Cells.Item(45, 3).Value = "https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana"
$ws.Cells.Item(45, 4).NumberFormat = "@"
$ws.Cells.Item(45, 4).Value = "0.5360"
$ws.Cells.Item(45, 5).Value = "  -1.47%  "

$ws.Cells.Item(46, 2).Value = "RenderToken"
$ws.Cells.Item(46, 3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(46, 4).NumberFormat = "@"
$ws.Cells.Item(46, 4).Value = "2.188"
$ws.Cells.Item(46, 5).Value = "  -1.96%  "

$ws.Cells.Item(47, 2).Value = "Cronos"
$ws.Cells.Item(47, 3).Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Cells.Item(47, 4).NumberFormat = "@"
$ws.Cells.Item(47, 4).Value = "0.07023"
$ws.Cells.Item(47, 5).Value = "  -1.02%  "

$ws.Cells.Item(48, 2).Value = "NEARProtocol"
$ws.Cells.Item(48, 3).Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Cells.Item(48, 4).NumberFormat = "@"
$ws.Cells.Item(48, 4).Value = "1.851"
$ws.Cells.Item(48, 5).Value = "  -1.99%  "

$ws.Cells.Item(49, 2).Value = "MXToken"
$ws.Cells.Item(49, 3).Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Cells.Item(49, 4).NumberFormat = "@"
$ws.Cells.Item(49, 4).Value = "2.571"
$ws.Cells.Item(49, 5).Value = "  -0.69%  "

$ws.Cells.Item(50, 2).Value = "Quant"
$ws.Cells.Item(50, 3).Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Cells.Item(50, 4).NumberFormat = "@"
$ws.Cells.Item(50, 4).Value = "112.97"
$ws.Cells.Item(50, 5).Value = "  +0.79%  "

$ws.Cells.Item(51, 2).Value = "WEMIXToken"
$ws.Cells.Item(51, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(51, 4).NumberFormat = "@"
$ws.Cells.Item(51, 4).Value = "1.046"
$ws.Cells.Item(51, 5).Value = "  -2.37%  "
